# Applies updated crypto price/volume figures (and restores the Kaspa/BitcoinSV
# row order) per the Oct 20 2023 GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.470.12"
$ws.Range("E2").Value = "  +2.97%  "

$ws.Range("D3").Value = "1.602.97"
$ws.Range("E3").Value = "  +2.48%  "

$ws.Range("D4").Value = "`'0.999"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "`'212.28"
$ws.Range("E5").Value = "  +0.86%  "

$ws.Range("E6").Value = "  +6.67%  "

$ws.Range("D7").Value = "`'0.999"
$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").Value = "`'26.47"
$ws.Range("E8").Value = "  +6.05%  "

$ws.Range("D9").Value = "`'43.39"
$ws.Range("E9").Value = "  -1.11%  "

$ws.Range("E10").Value = "  +2.39%  "

$ws.Range("D11").Value = "`'0.0599"
$ws.Range("E11").Value = "  +2.31%  "

$ws.Range("E12").Value = "  +1.65%  "

$ws.Range("D13").Value = "1.833.04"
$ws.Range("E13").Value = "  +2.49%  "

$ws.Range("D14").Value = "1.609.08"
$ws.Range("E14").Value = "  +2.97%  "

$ws.Range("D15").Value = "29.504.71"
$ws.Range("E15").Value = "  +2.96%  "

$ws.Range("E16").Value = "  +3.48%  "

$ws.Range("E17").Value = "  +1.44%  "

$ws.Range("D18").Value = "`'63.25"
$ws.Range("E18").Value = "  +3.06%  "

$ws.Range("D19").Value = "`'240.65"
$ws.Range("E19").Value = "  +4.03%  "

$ws.Range("E20").Value = "  +3.70%  "

$ws.Range("E21").Value = "  +2.22%  "

$ws.Range("E22").Value = "  -0.03%  "

$ws.Range("E23").Value = "  +1.84%  "

$ws.Range("D24").Value = "`'9.14"
$ws.Range("E24").Value = "  +1.71%  "

$ws.Range("D25").Value = "`'2.09"
$ws.Range("E25").Value = "  -2.19%  "

$ws.Range("E26").Value = "  +2.42%  "

$ws.Range("D27").Value = "`'15.28"
$ws.Range("E27").Value = "  +3.35%  "

$ws.Range("E28").Value = "  +5.11%  "

$ws.Range("E29").Value = "  +2.38%  "

$ws.Range("E30").Value = "  -0.04%  "

$ws.Range("E31").Value = "  +2.55%  "

$ws.Range("E32").Value = "  -0.17%  "

$ws.Range("E33").Value = "  +1.75%  "

$ws.Range("E34").Value = "  +3.91%  "

$ws.Range("D35").Value = "1.414.83"
$ws.Range("E35").Value = "  +1.90%  "

$ws.Range("E36").Value = "  -1.03%  "

$ws.Range("E37").Value = "  +2.92%  "

$ws.Range("D38").Value = "`'2.82"
$ws.Range("E38").Value = "  +6.29%  "

$ws.Range("E39").Value = "  +0.20%  "

$ws.Range("E40").Value = "  +2.14%  "

$ws.Range("E41").Value = "  +3.31%  "

$ws.Range("E42").Value = "  +0.97%  "

$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "`'0.0484"
$ws.Range("E43").Value = "  +5.52%  "

$ws.Range("B44").Value = "BitcoinSV"
$ws.Range("C44").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D44").Value = "`'53.31"
$ws.Range("E44").Value = "  +23.84%  "

$ws.Range("E45").Value = "  -0.05%  "

$ws.Range("E46").Value = "  +2.28%  "

$ws.Range("E47").Value = "  +2.70%  "

$ws.Range("E48").Value = "  +0.18%  "

$ws.Range("D49").Value = "1.743.94"
$ws.Range("E49").Value = "  +2.57%  "

$ws.Range("D50").Value = "`'86.54"
$ws.Range("E50").Value = "  +1.44%  "

$ws.Range("D51").Value = "`'0.844"
$ws.Range("E51").Value = "  -2.76%  "
